# The author's edit:
#  1. Changed the "queries" sheet's B2 cell text from
#     "select * from projectStructureDetails" to "select * from projectLevelDetails"
#  2. Updated the selection/active-cell on "projectLevelDetails" sheet to G10
#  3. Updated the selection/active-cell on "queries" sheet to B8
#  4. Made "queries" the active (selected) sheet/tab instead of "testCasesTestNG"

$wb = $excel.ActiveWorkbook

# --- Update the selection on projectLevelDetails (first sheet) ---
$wsProjectLevelDetails = $wb.Worksheets.Item("projectLevelDetails")
$wsProjectLevelDetails.Activate()
$wsProjectLevelDetails.Range("G10").Select() | Out-Null

# --- Update the query text and selection on the queries sheet, and make it active ---
$wsQueries = $wb.Worksheets.Item("queries")
$wsQueries.Range("B2").Value = "select * from projectLevelDetails"
$wsQueries.Activate()
$wsQueries.Range("B8").Select() | Out-Null
